$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "rank" in B1 to "cost" (moved buy cost of upgrade card to the center)
$ws.Range("B1").Value = "cost"

# Update selection to B2 as in diff
$ws.Range("B2").Select()
